$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 7642.5713
$ws.Range("I69").Value = 6699.6
$ws.Range("J69").Value = 10000
$ws.Range("K69").Value = 20098.8
$ws.Range("L69").Value = 30000
$ws.Range("M69").Value = -19224.8
$ws.Range("N69").Value = -31748

$ws.Range("H72").Value = 7642.5713
$ws.Range("I72").Value = 6699.6
$ws.Range("J72").Value = 10000
$ws.Range("K72").Value = 60296.4
$ws.Range("L72").Value = 90000
$ws.Range("M72").Value = -55928.4
$ws.Range("N72").Value = -98736

$ws.Range("H98").Value = 1012.73334
$ws.Range("I98").Value = 630.8461
$ws.Range("K98").Value = 630.8461
$ws.Range("M98").Value = 867.1539

$ws.Range("H113").Value = 7899.85
$ws.Range("I113").Value = 7830.8667
$ws.Range("J113").Value = 8106.8
$ws.Range("K113").Value = 7830.8667
$ws.Range("L113").Value = 8106.8
$ws.Range("M113").Value = -4576.8667
$ws.Range("N113").Value = -14614.8

$ws.Range("H122").Value = 1012.73334
$ws.Range("I122").Value = 630.8461
$ws.Range("K122").Value = 1892.5383
$ws.Range("M122").Value = 557.4617000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 448.875
$ws.Range("I4").Value = 391.64285
$ws.Range("K4").Value = 391.64285
$ws.Range("M4").Value = -275.64285

$ws.Range("H32").Value = 5755.0835
$ws.Range("I32").Value = 2453.492
$ws.Range("K32").Value = 2453.492
$ws.Range("M32").Value = -2166.492

$ws.Range("H45").Value = 3200
$ws.Range("I45").Value = 1666.6666
$ws.Range("K45").Value = 1666.6666
$ws.Range("M45").Value = -1289.6666

$ws.Range("H61").Value = 2636.889
$ws.Range("I61").Value = 1842.9333
$ws.Range("K61").Value = 1842.9333
$ws.Range("M61").Value = -1630.9333

$ws.Range("H102").Value = 3055
$ws.Range("I102").Value = 1350.6
$ws.Range("K102").Value = 1350.6
$ws.Range("M102").Value = 271.4000000000001

$ws.Range("H110").Value = 2500
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H122").Value = 1743.5454
$ws.Range("I122").Value = 1687.9
$ws.Range("K122").Value = 5063.700000000001
$ws.Range("M122").Value = -2613.700000000001

$ws.Range("H136").Value = 2636.889
$ws.Range("I136").Value = 1842.9333
$ws.Range("K136").Value = 5528.7999
$ws.Range("M136").Value = -2978.7999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1576.5
$ws.Range("I7").Value = 1103
$ws.Range("J7").Value = 2050
$ws.Range("K7").Value = 1103
$ws.Range("L7").Value = 2050
$ws.Range("M7").Value = -990
$ws.Range("N7").Value = -2276

$ws.Range("H20").Value = 447
$ws.Range("I20").Value = 309
$ws.Range("K20").Value = 309
$ws.Range("M20").Value = -62

$ws.Range("H86").Value = 1626.0714
$ws.Range("I86").Value = 1688.8334
$ws.Range("J86").Value = 1249.5
$ws.Range("K86").Value = 1688.8334
$ws.Range("L86").Value = 1249.5
$ws.Range("M86").Value = -565.8334
$ws.Range("N86").Value = -3495.5

$ws.Range("H89").Value = 1626.0714
$ws.Range("I89").Value = 1688.8334
$ws.Range("J89").Value = 1249.5
$ws.Range("K89").Value = 8444.166999999999
$ws.Range("L89").Value = 6247.5
$ws.Range("M89").Value = -2828.166999999999
$ws.Range("N89").Value = -17479.5

$ws.Range("H122").Value = 50000
$ws.Range("J122").Value = 50000
$ws.Range("L122").Value = 50000
$ws.Range("N122").Value = -59800

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 176.64285
$ws.Range("J7").Value = 135.5
$ws.Range("L7").Value = 135.5
$ws.Range("N7").Value = -361.5

$ws.Range("H31").Value = 7444.8
$ws.Range("J31").Value = 8310.272000000001
$ws.Range("L31").Value = 8310.272000000001
$ws.Range("N31").Value = -8900.272000000001

$ws.Range("H34").Value = 7444.8
$ws.Range("J34").Value = 8310.272000000001
$ws.Range("L34").Value = 8310.272000000001
$ws.Range("N34").Value = -8714.272000000001

$ws.Range("H99").Value = 1824
$ws.Range("I99").Value = 1728
$ws.Range("K99").Value = 1728
$ws.Range("M99").Value = -230

$ws.Range("H105").Value = 1990.6666
$ws.Range("I105").Value = 1916.9445
$ws.Range("J105").Value = 2433
$ws.Range("K105").Value = 1916.9445
$ws.Range("L105").Value = 2433
$ws.Range("M105").Value = -169.9445000000001
$ws.Range("N105").Value = -5927

$ws.Range("H126").Value = 1824
$ws.Range("I126").Value = 1728
$ws.Range("K126").Value = 5184
$ws.Range("M126").Value = -2714

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 384662.53
$ws.Range("J2").Value = 125
$ws.Range("L2").Value = 750
$ws.Range("N2").Value = -976

$ws.Range("H37").Value = 74249.25
$ws.Range("J37").Value = 74249.25
$ws.Range("L37").Value = 222747.75
$ws.Range("N37").Value = -222971.75

$ws.Range("H113").Value = 934.4
$ws.Range("J113").Value = 968.25
$ws.Range("L113").Value = 2904.75
$ws.Range("N113").Value = -7244.75

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

$ws.Range("H131").Value = 1174268.1
$ws.Range("I131").Value = 124231
$ws.Range("J131").Value = 1324273.4
$ws.Range("K131").Value = 372693
$ws.Range("L131").Value = 3972820.2
$ws.Range("M131").Value = -367653
$ws.Range("N131").Value = -3982900.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 103
$ws.Range("I2").Value = 119.166664
$ws.Range("K2").Value = 119.166664
$ws.Range("M2").Value = -6.166663999999997

$ws.Range("H13").Value = 235.42857
$ws.Range("I13").Value = 99.333336
$ws.Range("J13").Value = 337.5
$ws.Range("K13").Value = 99.333336
$ws.Range("L13").Value = 337.5
$ws.Range("M13").Value = 39.666664
$ws.Range("N13").Value = -615.5

$ws.Range("H63").Value = 50114
$ws.Range("J63").Value = 50114
$ws.Range("L63").Value = 50114
$ws.Range("N63").Value = -51486

$ws.Range("H66").Value = 50114
$ws.Range("J66").Value = 50114
$ws.Range("L66").Value = 150342
$ws.Range("N66").Value = -157206

$ws.Range("H102").Value = 2697.111
$ws.Range("I102").Value = 1580.6154
$ws.Range("J102").Value = 5600
$ws.Range("K102").Value = 1580.6154
$ws.Range("L102").Value = 5600
$ws.Range("M102").Value = 41.38460000000009
$ws.Range("N102").Value = -8844

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3971.9524
$ws.Range("I40").Value = 4062.9333
$ws.Range("J40").Value = 3744.5
$ws.Range("K40").Value = 4062.9333
$ws.Range("L40").Value = 3744.5
$ws.Range("M40").Value = -3926.9333
$ws.Range("N40").Value = -4016.5

$ws.Range("H122").Value = 4076.8572
$ws.Range("J122").Value = 4759.4
$ws.Range("L122").Value = 14278.2
$ws.Range("N122").Value = -19178.2

$ws.Range("H132").Value = 2667.6667
$ws.Range("J132").Value = 4333.1665
$ws.Range("L132").Value = 12999.4995
$ws.Range("N132").Value = -18059.4995

$ws.Range("H136").Value = 2295.5264
$ws.Range("I136").Value = 1059.6666
$ws.Range("K136").Value = 3178.9998
$ws.Range("M136").Value = -628.9998000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 3818.8
$ws.Range("I14").Value = 4748.5
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 4748.5
$ws.Range("L14").Value = 100
$ws.Range("M14").Value = -4580.5
$ws.Range("N14").Value = -436

$ws.Range("H34").Value = 49999
$ws.Range("I34").Value = 49999
$ws.Range("K34").Value = 49999
$ws.Range("M34").Value = -49796

$ws.Range("H107").Value = 1193.881
$ws.Range("I107").Value = 1354
$ws.Range("K107").Value = 4062
$ws.Range("M107").Value = -2142
